{"js": "// Fix extra whitespace: the run \"   \u05d1\u05d7\u05d5\u05e8\u05e3 \u05d4\u05d8\u05de\" (3 leading spaces) becomes\n// \" \u05d1\u05d7\u05d5\u05e8\u05e3 \u05d4\u05d8\u05de\" (1 leading space).\nconst searchResults = context.document.body.search(\"   \u05d1\u05d7\u05d5\u05e8\u05e3 \u05d4\u05d8\u05de\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\n\nsearchResults.items[0].insertText(\" \u05d1\u05d7\u05d5\u05e8\u05e3 \u05d4\u05d8\u05de\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix extra whitespace: the run \"   \u05d1\u05d7\u05d5\u05e8\u05e3 \u05d4\u05d8\u05de\" (3 leading spaces) becomes\n# \" \u05d1\u05d7\u05d5\u05e8\u05e3 \u05d4\u05d8\u05de\" (1 leading space) - i.e. \"...\u05d0\u05dd \u05dc\u05d0.   \u05d1\u05d7\u05d5\u05e8\u05e3...\" -> \"...\u05d0\u05dd \u05dc\u05d0. \u05d1\u05d7\u05d5\u05e8\u05e3...\".\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Replacement.ClearFormatting()\n$range.Find.Text = \"   \u05d1\u05d7\u05d5\u05e8\u05e3 \u05d4\u05d8\u05de\"\n$range.Find.Replacement.Text = \" \u05d1\u05d7\u05d5\u05e8\u05e3 \u05d4\u05d8\u05de\"\n$range.Find.Forward = $true\n$range.Find.Wrap = 0\n$range.Find.Format = $false\n$range.Find.MatchCase = $true\n$range.Find.MatchWholeWord = $false\n$range.Find.MatchWildcards = $false\n\n$range.Find.Execute(\n    $range.Find.Text,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    0,\n    $false,\n    $range.Find.Replacement.Text,\n    1\n)\n"}
